$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

Set-CellText 'D2' '76.280.43'
Set-CellText 'E2' '  +0.47%  '
Set-CellText 'D3' '2.934.07'
Set-CellText 'E3' '  +3.29%  '
Set-CellText 'D4' '0.999'
Set-CellText 'E4' '  -0.02%  '
Set-CellText 'D5' '200.45'
Set-CellText 'E5' '  +4.50%  '
Set-CellText 'D6' '597.36'
Set-CellText 'E6' '  +0.02%  '
Set-CellText 'D7' '0.999'
Set-CellText 'E7' '  -0.05%  '
Set-CellText 'E8' '  +0.46%  '
Set-CellText 'D9' '0.198'
Set-CellText 'E9' '  +2.38%  '
Set-CellText 'D10' '2.930.13'
Set-CellText 'E10' '  +3.25%  '
Set-CellText 'D11' '0.448'
Set-CellText 'E11' '  +16.35%  '
Set-CellText 'E12' '  +0.71%  '
Set-CellText 'D13' '4.96'
Set-CellText 'D14' '3.463.28'
Set-CellText 'E14' '  +3.21%  '
Set-CellText 'D15' '76.110.96'
Set-CellText 'E15' '  +0.47%  '
Set-CellText 'D16' '27.98'
Set-CellText 'E16' '  +3.11%  '
Set-CellText 'E17' '  +0.36%  '
Set-CellText 'D18' '2.917.43'
Set-CellText 'E18' '  +2.70%  '
Set-CellText 'D19' '13.31'
Set-CellText 'E19' '  +7.40%  '
Set-CellText 'D20' '8.76'
Set-CellText 'E20' '  -3.40%  '
Set-CellText 'D21' '372.47'
Set-CellText 'E21' '  -2.22%  '
Set-CellText 'E22' '  -1.21%  '
Set-CellText 'E23' '  +4.51%  '
Set-CellText 'D24' '72.51'
Set-CellText 'E24' '  +1.64%  '
Set-CellText 'D25' '0.999'
Set-CellText 'E25' '  -0.03%  '
Set-CellText 'D27' '4.27'
Set-CellText 'E27' '  +1.50%  '
Set-CellText 'D28' '9.66'
Set-CellText 'E28' '  -1.50%  '
Set-CellText 'E29' '  +2.67%  '
Set-CellText 'D30' '0.999'
Set-CellText 'E30' '  +0.07%  '
Set-CellText 'E31' '  -2.83%  '
Set-CellText 'E32' '  +2.16%  '
Set-CellText 'D33' '498.26'
Set-CellText 'E33' '  -3.73%  '
Set-CellText 'D34' '1.84'
Set-CellText 'E34' '  +1.13%  '
Set-CellText 'D35' '0.999'
Set-CellText 'E35' '  +0.04%  '
Set-CellText 'D36' '164.39'
Set-CellText 'E36' '  +0.32%  '
Set-CellText 'D37' '20.12'
Set-CellText 'E37' '  +0.67%  '
Set-CellText 'D38' '0.389'
Set-CellText 'E38' '  +13.13%  '
Set-CellText 'E39' '  +23.86%  '
Set-CellText 'E41' '  -6.55%  '
Set-CellText 'E42' '  -0.01%  '
Set-CellText 'D43' '178.64'
Set-CellText 'E43' '  -3.06%  '
Set-CellText 'D44' '4.93'
Set-CellText 'E44' '  -2.70%  '
Set-CellText 'E45' '  -1.83%  '
Set-CellText 'D46' '40.16'
Set-CellText 'E46' '  +0.38%  '
Set-CellText 'E47' '  -1.69%  '
Set-CellText 'B48' 'ARBITRUM'
Set-CellText 'C48' 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-CellText 'D48' '0.584'
Set-CellText 'E48' '  +1.64%  '
Set-CellText 'B49' 'Filecoin'
Set-CellText 'C49' 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-CellText 'D49' '3.88'
Set-CellText 'E49' '  +3.35%  '
Set-CellText 'D50' '2.32'
Set-CellText 'E50' '  -1.91%  '
Set-CellText 'D51' '22.85'
Set-CellText 'E51' '  +7.76%  '
